$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '34.549.64'
$ws.Range('E2').Value = '  -0.17%  '

# Row 3
$ws.Range('D3').Value = '1.810.80'
$ws.Range('E3').Value = '  -0.34%  '

# Row 4
$ws.Range('E4').Value = '  +0.08%  '

# Row 5
$ws.Range('D5').Value = '''228.53'
$ws.Range('E5').Value = '  +0.06%  '

# Row 6
$ws.Range('D6').Value = '''0.606'
$ws.Range('E6').Value = '  +8.33%  '

# Row 7
$ws.Range('E7').Value = '  -0.01%  '

# Row 8
$ws.Range('D8').Value = '''36.63'
$ws.Range('E8').Value = '  +5.22%  '

# Row 9
$ws.Range('D9').Value = '''0.301'
$ws.Range('E9').Value = '  +0.23%  '

# Row 10
$ws.Range('D10').Value = '''0.0700'
$ws.Range('E10').Value = '  +0.72%  '

# Row 11
$ws.Range('D11').Value = '''0.0965'
$ws.Range('E11').Value = '  +1.31%  '

# Row 12
$ws.Range('D12').Value = '2.070.72'
$ws.Range('E12').Value = '  -0.33%  '

# Row 13
$ws.Range('D13').Value = '''11.48'
$ws.Range('E13').Value = '  +1.18%  '

# Row 14
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.825.14'
$ws.Range('E14').Value = '  +0.31%  '

# Row 15
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').Value = '''0.655'
$ws.Range('E15').Value = '  +1.73%  '

# Row 16
$ws.Range('E16').Value = '  +3.99%  '

# Row 17
$ws.Range('D17').Value = '34.549.63'
$ws.Range('E17').Value = '  -0.26%  '

# Row 18
$ws.Range('D18').Value = '''70.15'
$ws.Range('E18').Value = '  +1.38%  '

# Row 19
$ws.Range('D19').Value = '''246.95'
$ws.Range('E19').Value = '  -0.28%  '

# Row 20
$ws.Range('D20').Value = '0.0₃0799'
$ws.Range('E20').Value = '  -0.46%  '

# Row 21
$ws.Range('D21').Value = '''11.63'
$ws.Range('E21').Value = '  +0.69%  '

# Row 22
$ws.Range('E22').Value = '  -0.01%  '

# Row 23
$ws.Range('D23').Value = '''4.22'
$ws.Range('E23').Value = '  +0.16%  '

# Row 24
$ws.Range('D24').Value = '''2.24'
$ws.Range('E24').Value = '  +7.09%  '

# Row 25
$ws.Range('D25').Value = '''172.99'
$ws.Range('E25').Value = '  +0.42%  '

# Row 26
$ws.Range('E26').Value = '  +7.39%  '

# Row 27
$ws.Range('D27').Value = '''17.34'
$ws.Range('E27').Value = '  +3.40%  '

# Row 28
$ws.Range('D28').Value = '''0.123'
$ws.Range('E28').Value = '  +5.08%  '

# Row 29
$ws.Range('E29').Value = '  -0.04%  '

# Row 30
$ws.Range('D30').Value = '''4.03'
$ws.Range('E30').Value = '  -0.07%  '

# Row 31
$ws.Range('D31').Value = '''3.86'
$ws.Range('E31').Value = '  +0.15%  '

# Row 32
$ws.Range('D32').Value = '''0.0531'
$ws.Range('E32').Value = '  -0.35%  '

# Row 33
$ws.Range('D33').Value = '''1.24'

# Row 34
$ws.Range('E34').Value = '  -1.70%  '

# Row 35
$ws.Range('D35').Value = '1.403.47'
$ws.Range('E35').Value = '  -1.18%  '

# Row 36
$ws.Range('D36').Value = '''0.675'
$ws.Range('E36').Value = '  -0.72%  '

# Row 37
$ws.Range('D37').Value = '''2.45'
$ws.Range('E37').Value = '  -6.18%  '

# Row 38
$ws.Range('E38').Value = '  +0.40%  '

# Row 39
$ws.Range('E39').Value = '  -0.97%  '

# Row 40
$ws.Range('D40').Value = '''0.969'
$ws.Range('E40').Value = '  +0.49%  '

# Row 41
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').Value = '''83.00'
$ws.Range('E41').Value = '  -3.68%  '

# Row 42
$ws.Range('B42').Value = 'MXToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D42').Value = '''2.84'
$ws.Range('E42').Value = '  -0.66%  '

# Row 43
$ws.Range('E43').Value = '  +0.48%  '

# Row 44
$ws.Range('E44').Value = '  +7.51%  '

# Row 45
$ws.Range('E45').Value = '  -2.53%  '

# Row 46
$ws.Range('D46').Value = '''6.06'
$ws.Range('E46').Value = '  -0.92%  '

# Row 47
$ws.Range('B47').Value = 'RocketPoolETH'
$ws.Range('C47').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D47').Value = '1.972.02'
$ws.Range('E47').Value = '  -0.35%  '

# Row 48
$ws.Range('B48').Value = 'Kaspa'
$ws.Range('C48').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D48').Value = '''0.0495'
$ws.Range('E48').Value = '  -5.60%  '

# Row 49
$ws.Range('D49').Value = '''104.63'
$ws.Range('E49').Value = '  -1.46%  '

# Row 50
$ws.Range('E50').Value = '  -0.01%  '

# Row 51
$ws.Range('D51').Value = '0.0₆0127'
$ws.Range('E51').Value = '  -3.46%  '

